$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 811.5714
$ws.Range("I15").Value = 811.5714
$ws.Range("K15").Value = 2434.7142
$ws.Range("M15").Value = -2265.7142
$ws.Range("H19").Value = 696.1818
$ws.Range("I19").Value = 100.666664
$ws.Range("K19").Value = 100.666664
$ws.Range("M19").Value = 74.333336
$ws.Range("H28").Value = 606.95
$ws.Range("I28").Value = 349.33334
$ws.Range("K28").Value = 349.33334
$ws.Range("M28").Value = 135.66666
$ws.Range("H40").Value = 3229.8948
$ws.Range("I40").Value = 3377.875
$ws.Range("J40").Value = 2440.6667
$ws.Range("K40").Value = 3377.875
$ws.Range("L40").Value = 2440.6667
$ws.Range("M40").Value = -3202.875
$ws.Range("N40").Value = -2790.6667
$ws.Range("H61").Value = 85.333336
$ws.Range("I61").Value = 85.333336
$ws.Range("K61").Value = 256.000008
$ws.Range("M61").Value = -84.00000799999998
$ws.Range("H62").Value = 6508
$ws.Range("I62").Value = 6508
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 6508
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -5884
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 6508
$ws.Range("I65").Value = 6508
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 32540
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -29420
$ws.Range("N65").ClearContents()
$ws.Range("H70").Value = 7741.3687
$ws.Range("J70").Value = 7811.343
$ws.Range("L70").Value = 23434.029
$ws.Range("N70").Value = -23974.029
$ws.Range("H73").Value = 7741.3687
$ws.Range("J73").Value = 7811.343
$ws.Range("L73").Value = 23434.029
$ws.Range("N73").Value = -25306.029
$ws.Range("H96").Value = 235.84616
$ws.Range("I96").Value = 262.44446
$ws.Range("J96").Value = 176
$ws.Range("K96").Value = 787.33338
$ws.Range("L96").Value = 528
$ws.Range("M96").Value = 585.66662
$ws.Range("N96").Value = -3274
$ws.Range("H98").Value = 1134.2903
$ws.Range("I98").Value = 1145.434
$ws.Range("K98").Value = 1145.434
$ws.Range("M98").Value = 352.566
$ws.Range("H103").Value = 649.2222
$ws.Range("I103").Value = 763.1667
$ws.Range("J103").Value = 421.33334
$ws.Range("K103").Value = 2289.5001
$ws.Range("L103").Value = 1264.00002
$ws.Range("M103").Value = -1703.5001
$ws.Range("N103").Value = -2436.00002
$ws.Range("H116").Value = 5625.593
$ws.Range("I116").Value = 4212.467
$ws.Range("K116").Value = 4212.467
$ws.Range("M116").Value = -770.4669999999996
$ws.Range("H122").Value = 1134.2903
$ws.Range("I122").Value = 1145.434
$ws.Range("K122").Value = 3436.302
$ws.Range("M122").Value = -986.3019999999997
$ws.Range("H132").Value = 22728930
$ws.Range("I132").Value = 24391974
$ws.Range("K132").Value = 73175922
$ws.Range("M132").Value = -73173392
$ws.Range("H137").Value = 51745.473
$ws.Range("I137").Value = 70233.30499999999
$ws.Range("K137").Value = 210699.915
$ws.Range("M137").Value = -208149.915
$ws.Range("H141").Value = 9555.875
$ws.Range("I141").Value = 9859.6
$ws.Range("K141").Value = 29578.8
$ws.Range("M141").Value = -24398.8

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H41").Value = 1874.2
$ws.Range("I41").Value = 1776
$ws.Range("J41").Value = 1939.6666
$ws.Range("K41").Value = 1776
$ws.Range("L41").Value = 1939.6666
$ws.Range("M41").Value = -1362
$ws.Range("N41").Value = -2767.6666
$ws.Range("H45").Value = 6692992
$ws.Range("I45").Value = 8549324
$ws.Range("K45").Value = 8549324
$ws.Range("M45").Value = -8548947
$ws.Range("H61").Value = 4069.9092
$ws.Range("I61").Value = 3982.15
$ws.Range("J61").Value = 4947.5
$ws.Range("K61").Value = 3982.15
$ws.Range("L61").Value = 4947.5
$ws.Range("M61").Value = -3770.15
$ws.Range("N61").Value = -5371.5
$ws.Range("H63").Value = 4754.8945
$ws.Range("I63").Value = 2313.9092
$ws.Range("J63").Value = 8111.25
$ws.Range("K63").Value = 2313.9092
$ws.Range("L63").Value = 8111.25
$ws.Range("M63").Value = -1627.9092
$ws.Range("N63").Value = -9483.25
$ws.Range("H66").Value = 4754.8945
$ws.Range("I66").Value = 2313.9092
$ws.Range("J66").Value = 8111.25
$ws.Range("K66").Value = 11569.546
$ws.Range("L66").Value = 40556.25
$ws.Range("M66").Value = -8137.546
$ws.Range("N66").Value = -47420.25
$ws.Range("H74").Value = 30724.162
$ws.Range("I74").Value = 1341.826
$ws.Range("K74").Value = 1341.826
$ws.Range("M74").Value = -467.826
$ws.Range("H77").Value = 30724.162
$ws.Range("I77").Value = 1341.826
$ws.Range("K77").Value = 6709.13
$ws.Range("M77").Value = -2341.13
$ws.Range("H97").Value = 1294721.1
$ws.Range("I97").Value = 1903553
$ws.Range("J97").Value = 953.25
$ws.Range("K97").Value = 1903553
$ws.Range("L97").Value = 953.25
$ws.Range("M97").Value = -1903057
$ws.Range("N97").Value = -1945.25
$ws.Range("H102").Value = 4389387
$ws.Range("I102").Value = 5955030.5
$ws.Range("J102").Value = 5586
$ws.Range("K102").Value = 5955030.5
$ws.Range("L102").Value = 5586
$ws.Range("M102").Value = -5953408.5
$ws.Range("N102").Value = -8830
$ws.Range("H122").Value = 466422.38
$ws.Range("I122").Value = 2259.4138
$ws.Range("K122").Value = 6778.241399999999
$ws.Range("M122").Value = -4328.241399999999
$ws.Range("H132").Value = 2281.6584
$ws.Range("I132").Value = 1803.5264
$ws.Range("J132").Value = 8338
$ws.Range("K132").Value = 5410.5792
$ws.Range("L132").Value = 25014
$ws.Range("M132").Value = -2880.5792
$ws.Range("N132").Value = -30074
$ws.Range("H136").Value = 4069.9092
$ws.Range("I136").Value = 3982.15
$ws.Range("J136").Value = 4947.5
$ws.Range("K136").Value = 11946.45
$ws.Range("L136").Value = 14842.5
$ws.Range("M136").Value = -9396.450000000001
$ws.Range("N136").Value = -19942.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I23").Value = 1006
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 1006
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -723
$ws.Range("N23").ClearContents()
$ws.Range("H40").Value = 45000
$ws.Range("J40").Value = 45000
$ws.Range("L40").Value = 45000
$ws.Range("N40").Value = -45530
$ws.Range("H94").Value = 3578949.2
$ws.Range("I94").Value = 5000563.5
$ws.Range("J94").Value = 24913.75
$ws.Range("K94").Value = 5000563.5
$ws.Range("L94").Value = 24913.75
$ws.Range("M94").Value = -5000112.5
$ws.Range("N94").Value = -25815.75
$ws.Range("H107").Value = 3404007
$ws.Range("I107").Value = 4203656
$ws.Range("K107").Value = 4203656
$ws.Range("M107").Value = -4201736

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1783.6
$ws.Range("I16").Value = 1265.4286
$ws.Range("K16").Value = 1265.4286
$ws.Range("M16").Value = -978.4286
$ws.Range("H31").Value = 22952.225
$ws.Range("I31").Value = 2732.7
$ws.Range("J31").Value = 28136.719
$ws.Range("K31").Value = 2732.7
$ws.Range("L31").Value = 28136.719
$ws.Range("M31").Value = -2437.7
$ws.Range("N31").Value = -28726.719
$ws.Range("H34").Value = 22952.225
$ws.Range("I34").Value = 2732.7
$ws.Range("J34").Value = 28136.719
$ws.Range("K34").Value = 2732.7
$ws.Range("L34").Value = 28136.719
$ws.Range("M34").Value = -2530.7
$ws.Range("N34").Value = -28540.719
$ws.Range("H38").Value = 15000
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()
$ws.Range("H46").Value = 15000
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("H58").Value = 6238.148
$ws.Range("I58").Value = 7026.8423
$ws.Range("J58").Value = 4365
$ws.Range("K58").Value = 7026.8423
$ws.Range("L58").Value = 4365
$ws.Range("M58").Value = -6823.8423
$ws.Range("N58").Value = -4771
$ws.Range("H62").Value = 4747.5
$ws.Range("J62").Value = 4747.5
$ws.Range("L62").Value = 4747.5
$ws.Range("N62").Value = -5995.5
$ws.Range("H65").Value = 4747.5
$ws.Range("J65").Value = 4747.5
$ws.Range("L65").Value = 23737.5
$ws.Range("N65").Value = -29977.5
$ws.Range("H105").Value = 1021.2
$ws.Range("I105").Value = 1021.2
$ws.Range("K105").Value = 1021.2
$ws.Range("M105").Value = 725.8
$ws.Range("H113").Value = 1783.6
$ws.Range("I113").Value = 1265.4286
$ws.Range("K113").Value = 1265.4286
$ws.Range("M113").Value = 904.5714
$ws.Range("H132").Value = 62811.062
$ws.Range("I132").Value = 44836.25
$ws.Range("J132").Value = 116735.5
$ws.Range("K132").Value = 134508.75
$ws.Range("L132").Value = 350206.5
$ws.Range("M132").Value = -131978.75
$ws.Range("N132").Value = -355266.5
$ws.Range("H136").Value = 6238.148
$ws.Range("I136").Value = 7026.8423
$ws.Range("J136").Value = 4365
$ws.Range("K136").Value = 21080.5269
$ws.Range("L136").Value = 13095
$ws.Range("M136").Value = -18530.5269
$ws.Range("N136").Value = -18195
$ws.Range("H141").Value = 673441.7
$ws.Range("J141").Value = 673441.7
$ws.Range("L141").Value = 673441.7
$ws.Range("N141").Value = -683801.7

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 7215.7144
$ws.Range("I33").Value = 69.3
$ws.Range("J33").Value = 25081.75
$ws.Range("K33").Value = 415.8
$ws.Range("L33").Value = 150490.5
$ws.Range("M33").Value = -132.8
$ws.Range("N33").Value = -151056.5
$ws.Range("H107").Value = 2418.2856
$ws.Range("I107").Value = 2732.25
$ws.Range("J107").Value = 1999.6666
$ws.Range("K107").Value = 8196.75
$ws.Range("L107").Value = 5998.9998
$ws.Range("M107").Value = -6276.75
$ws.Range("N107").Value = -9838.9998
$ws.Range("H129").Value = 1721.5
$ws.Range("I129").Value = 943.6667
$ws.Range("J129").Value = 2499.3333
$ws.Range("K129").Value = 2831.0001
$ws.Range("L129").Value = 7497.999899999999
$ws.Range("M129").Value = 2168.9999
$ws.Range("N129").Value = -17497.9999
$ws.Range("H131").Value = 9923606
$ws.Range("I131").Value = 8334114.5
$ws.Range("J131").Value = 10420323
$ws.Range("K131").Value = 25002343.5
$ws.Range("L131").Value = 31260969
$ws.Range("M131").Value = -24997303.5
$ws.Range("N131").Value = -31271049
$ws.Range("H134").Value = 1828.8462
$ws.Range("I134").Value = 1828.8462
$ws.Range("K134").Value = 5486.5386
$ws.Range("M134").Value = -416.5385999999999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H38").Value = 39999
$ws.Range("J38").Value = 39999
$ws.Range("L38").Value = 39999
$ws.Range("N38").Value = -40925
$ws.Range("H107").Value = 226.4375
$ws.Range("I107").Value = 108.5
$ws.Range("K107").Value = 108.5
$ws.Range("M107").Value = 1811.5
$ws.Range("H122").Value = 319605.06
$ws.Range("I122").Value = 424976.34
$ws.Range("K122").Value = 1274929.02
$ws.Range("M122").Value = -1272479.02
$ws.Range("H135").Value = 73299.664
$ws.Range("J135").Value = 70499.5
$ws.Range("L135").Value = 70499.5
$ws.Range("N135").Value = -80639.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7014.0454
$ws.Range("I7").Value = 3329.2222
$ws.Range("J7").Value = 9565.076999999999
$ws.Range("K7").Value = 3329.2222
$ws.Range("L7").Value = 9565.076999999999
$ws.Range("M7").Value = -3217.2222
$ws.Range("N7").Value = -9789.076999999999
$ws.Range("H16").Value = 333.85715
$ws.Range("I16").Value = 333.85715
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 333.85715
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -163.85715
$ws.Range("N16").ClearContents()
$ws.Range("H96").Value = 59596
$ws.Range("J96").Value = 59596
$ws.Range("L96").Value = 59596
$ws.Range("N96").Value = -65088
$ws.Range("H122").Value = 6132.8184
$ws.Range("I122").Value = 3582.5715
$ws.Range("K122").Value = 10747.7145
$ws.Range("M122").Value = -8297.7145
$ws.Range("H126").Value = 7014.0454
$ws.Range("I126").Value = 3329.2222
$ws.Range("J126").Value = 9565.076999999999
$ws.Range("K126").Value = 9987.6666
$ws.Range("L126").Value = 28695.231
$ws.Range("M126").Value = -7517.6666
$ws.Range("N126").Value = -33635.231

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 5319.2
$ws.Range("J29").Value = 5319.2
$ws.Range("L29").Value = 5319.2
$ws.Range("N29").Value = -5899.2
$ws.Range("H81").Value = 7939772
$ws.Range("I81").Value = 13890435
$ws.Range("K81").Value = 27780870
$ws.Range("M81").Value = -27779809
$ws.Range("H84").Value = 7939772
$ws.Range("I84").Value = 13890435
$ws.Range("K84").Value = 138904350
$ws.Range("M84").Value = -138899046
$ws.Range("H122").Value = 6016.3335
$ws.Range("I122").Value = 4025.75
$ws.Range("K122").Value = 12077.25
$ws.Range("M122").Value = -9627.25
$ws.Range("H126").Value = 1208.75
$ws.Range("I126").Value = 1215.15
$ws.Range("J126").Value = 1200.75
$ws.Range("K126").Value = 3645.45
$ws.Range("L126").Value = 3602.25
$ws.Range("M126").Value = -1175.45
$ws.Range("N126").Value = -8542.25
$ws.Range("H132").Value = 13663926
$ws.Range("I132").Value = 16132587
$ws.Range("K132").Value = 48397761
$ws.Range("M132").Value = -48395231
$ws.Range("H136").Value = 3889.8386
$ws.Range("I136").Value = 3472.625
$ws.Range("K136").Value = 10417.875
$ws.Range("M136").Value = -7867.875
$ws.Range("H139").Value = 160992
$ws.Range("J139").Value = 160992
$ws.Range("L139").Value = 160992
$ws.Range("N139").Value = -171272
